# Weekly fruit/vegetable price update: insert this week's new records
# (Zapallo / Camote) at the top of the data block (row 680), pushing the
# existing historical rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 680.
$ws.Range("680:681").EntireRow.Insert()

# Row 680: "1a (guarda)"
$ws.Range("A680").Value = 8
$ws.Range("B680").Value = "Terminal La Palmera de La Serena"
$ws.Range("C680").Value = "Coquimbo"
$ws.Range("D680").Value = 44753
$ws.Range("E680").Value = 4
$ws.Range("F680").Value = 100112045
$ws.Range("G680").Value = "Zapallo"
$ws.Range("H680").Value = "Camote"
$ws.Range("I680").Value = "1a (guarda)"
$ws.Range("J680").Value = 1800
$ws.Range("K680").Value = 900
$ws.Range("L680").Value = 1000
$ws.Range("M680").Value = 950
$ws.Range("N680").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O680").Value = "Región de O'Higgins"
$ws.Range("P680").Value = 950
$ws.Range("Q680").Value = 1
$ws.Range("R680").Value = "Hortaliza"

# Row 681: "2a (guarda)"
$ws.Range("A681").Value = 8
$ws.Range("B681").Value = "Terminal La Palmera de La Serena"
$ws.Range("C681").Value = "Coquimbo"
$ws.Range("D681").Value = 44753
$ws.Range("E681").Value = 4
$ws.Range("F681").Value = 100112045
$ws.Range("G681").Value = "Zapallo"
$ws.Range("H681").Value = "Camote"
$ws.Range("I681").Value = "2a (guarda)"
$ws.Range("J681").Value = 860
$ws.Range("K681").Value = 800
$ws.Range("L681").Value = 850
$ws.Range("M681").Value = 825
$ws.Range("N681").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O681").Value = "Región de O'Higgins"
$ws.Range("P681").Value = 825
$ws.Range("Q681").Value = 1
$ws.Range("R681").Value = "Hortaliza"
